# Commit: "I0 and IF added"
# Adds two new trailing columns, I ("I0") and J ("IF"), to the single
# worksheet, with a header in row 1 and numeric values for data rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should pick up the same look as the rest of row 1
# (bold font + border via style index 1) -- copy the formatting from the
# existing "IP" header cell (H1) onto I1:J1, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row I0 / IF values for rows 2-36 (row index -> [I, J]).
$data = @{
    2  = @(1, 3)
    3  = @(4, 5)
    4  = @(2, 4)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(6, 6)
    8  = @(1, 3)
    9  = @(3, 4)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(1, 4)
    14 = @(1, 4)
    15 = @(1, 3)
    16 = @(1, 6)
    17 = @(1, 4)
    18 = @(5, 7)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 4)
    22 = @(1, 4)
    23 = @(1, 4)
    24 = @(4, 6)
    25 = @(3, 5)
    26 = @(1, 1)
    27 = @(1, 6)
    28 = @(1, 4)
    29 = @(5, 6)
    30 = @(7, 7)
    31 = @(7, 9)
    32 = @(1, 4)
    33 = @(1, 3)
    34 = @(4, 6)
    35 = @(3, 4)
    36 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value  = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]   # column J
}
